$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 184.7
$ws.Range("I5").Value = 194.625
$ws.Range("J5").Value = 145
$ws.Range("K5").Value = 194.625
$ws.Range("L5").Value = 145
$ws.Range("M5").Value = -79.625
$ws.Range("N5").Value = -375
# Row 33
$ws.Range("H33").Value = 193.75
$ws.Range("I33").Value = 125
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 125
$ws.Range("L33").Value = 400
$ws.Range("M33").Value = 104
$ws.Range("N33").Value = -858
# Row 64
$ws.Range("H64").Value = 83336970
$ws.Range("I64").Value = 1000000000
$ws.Range("J64").Value = 3969.0908
$ws.Range("K64").Value = 1000000000
$ws.Range("L64").Value = 3969.0908
$ws.Range("M64").Value = -999999752
$ws.Range("N64").Value = -4465.0908
# Row 67
$ws.Range("H67").Value = 83336970
$ws.Range("I67").Value = 1000000000
$ws.Range("J67").Value = 3969.0908
$ws.Range("K67").Value = 1000000000
$ws.Range("L67").Value = 3969.0908
$ws.Range("M67").Value = -999999142
$ws.Range("N67").Value = -5685.0908
# Row 113
$ws.Range("H113").Value = 2148.25
$ws.Range("I113").Value = 1792.2727
$ws.Range("J113").Value = 2583.3333
$ws.Range("K113").Value = 1792.2727
$ws.Range("L113").Value = 2583.3333
$ws.Range("M113").Value = 1461.7273
$ws.Range("N113").Value = -9091.3333

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 755.8
$ws.Range("I2").Value = 623.2105
$ws.Range("J2").Value = 1175.6666
$ws.Range("K2").Value = 623.2105
$ws.Range("L2").Value = 1175.6666
$ws.Range("M2").Value = -510.2105
$ws.Range("N2").Value = -1401.6666
# Row 45
$ws.Range("H45").Value = 2340.8572
$ws.Range("I45").Value = 1937.3334
$ws.Range("J45").Value = 2643.5
$ws.Range("K45").Value = 1937.3334
$ws.Range("L45").Value = 2643.5
$ws.Range("M45").Value = -1560.3334
$ws.Range("N45").Value = -3397.5
# Row 56
$ws.Range("H56").Value = 11000
$ws.Range("J56").Value = 11000
$ws.Range("L56").Value = 11000
$ws.Range("N56").Value = -12484
# Row 104
$ws.Range("H104").Value = 42741.668
$ws.Range("J104").Value = 42741.668
$ws.Range("L104").Value = 42741.668
$ws.Range("N104").Value = -49729.668
# Row 116
$ws.Range("H116").Value = 755.8
$ws.Range("I116").Value = 623.2105
$ws.Range("J116").Value = 1175.6666
$ws.Range("K116").Value = 623.2105
$ws.Range("L116").Value = 1175.6666
$ws.Range("M116").Value = 1670.7895
$ws.Range("N116").Value = -5763.6666
# Row 122
$ws.Range("H122").Value = 1535.909
$ws.Range("I122").Value = 1052
$ws.Range("K122").Value = 3156
$ws.Range("M122").Value = -706
# Row 123
$ws.Range("H123").Value = 39000
$ws.Range("J123").Value = 39000
$ws.Range("L123").Value = 39000
$ws.Range("N123").Value = -48800
# Row 130
$ws.Range("H130").Value = 33408.777
$ws.Range("J130").Value = 33408.777
$ws.Range("L130").Value = 33408.777
$ws.Range("N130").Value = -43448.777
# Row 132
$ws.Range("H132").Value = 5488.353
$ws.Range("I132").Value = 3856.8684
$ws.Range("J132").Value = 10257.308
$ws.Range("K132").Value = 11570.6052
$ws.Range("L132").Value = 30771.924
$ws.Range("M132").Value = -9040.6052
$ws.Range("N132").Value = -35831.924

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 755.8
$ws.Range("I3").Value = 623.2105
$ws.Range("J3").Value = 1175.6666
$ws.Range("K3").Value = 623.2105
$ws.Range("L3").Value = 1175.6666
$ws.Range("M3").Value = -509.2105
$ws.Range("N3").Value = -1403.6666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 126088.25
$ws.Range("I22").Value = 250561
$ws.Range("J22").Value = 1615.5
$ws.Range("K22").Value = 250561
$ws.Range("L22").Value = 1615.5
$ws.Range("M22").Value = -250211
$ws.Range("N22").Value = -2315.5
# Row 132
$ws.Range("H132").Value = 27032134
$ws.Range("I132").Value = 47626084
$ws.Range("J132").Value = 2574.125
$ws.Range("K132").Value = 142878252
$ws.Range("L132").Value = 7722.375
$ws.Range("M132").Value = -142875722
$ws.Range("N132").Value = -12782.375
# Row 134
$ws.Range("H134").Value = 3429.6843
$ws.Range("I134").Value = 3687.2856
$ws.Range("J134").Value = 2708.4
$ws.Range("K134").Value = 11061.8568
$ws.Range("L134").Value = 8125.200000000001
$ws.Range("M134").Value = -8526.856800000001
$ws.Range("N134").Value = -13195.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 1559.8
$ws.Range("I55").Value = 599.6667
$ws.Range("J55").Value = 3000
$ws.Range("K55").Value = 1799.0001
$ws.Range("L55").Value = 9000
$ws.Range("M55").Value = -1622.0001
$ws.Range("N55").Value = -9354
# Row 58
$ws.Range("H58").Value = 2714.2856
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -2872
$ws.Range("N58").Value = -9256
# Row 62
$ws.Range("H62").Value = 2874.75
$ws.Range("J62").Value = 2874.75
$ws.Range("L62").Value = 8624.25
$ws.Range("N62").Value = -9996.25
# Row 65
$ws.Range("H65").Value = 2874.75
$ws.Range("J65").Value = 2874.75
$ws.Range("L65").Value = 25872.75
$ws.Range("N65").Value = -32736.75
# Row 87
$ws.Range("H87").Value = 2831.3333
$ws.Range("I87").Value = 2831.3333
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 8493.999899999999
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -7245.999899999999
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 2831.3333
$ws.Range("I90").Value = 2831.3333
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 25481.9997
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -19241.9997
$ws.Range("N90").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 1685.25
$ws.Range("I126").Value = 1232
$ws.Range("J126").Value = 2037.7778
$ws.Range("K126").Value = 3696
$ws.Range("L126").Value = 6113.3334
$ws.Range("M126").Value = -1226
$ws.Range("N126").Value = -11053.3334

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 11932.223
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 11932.223
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 11932.223
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -12522.223
# Row 27
$ws.Range("H27").Value = 11932.223
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 11932.223
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 11932.223
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -12146.223
# Row 55
$ws.Range("H55").Value = 249.17392
$ws.Range("I55").Value = 245.25
$ws.Range("J55").Value = 251.26666
$ws.Range("K55").Value = 245.25
$ws.Range("L55").Value = 251.26666
$ws.Range("M55").Value = -72.25
$ws.Range("N55").Value = -597.26666
# Row 122
$ws.Range("H122").Value = 43018.12
$ws.Range("I122").Value = 69012.87
$ws.Range("K122").Value = 207038.61
$ws.Range("M122").Value = -204588.61
# Row 132
$ws.Range("H132").Value = 29823.611
$ws.Range("I132").Value = 24967.166
$ws.Range("J132").Value = 32251.834
$ws.Range("K132").Value = 74901.49800000001
$ws.Range("L132").Value = 96755.50199999999
$ws.Range("M132").Value = -72371.49800000001
$ws.Range("N132").Value = -101815.502
